# Apply cryptocurrency price/volume updates (and the two-row reorderings
# for rows 38-43) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to keep a numeric-looking string as text
# (mirrors how Excel stores manually quote-prefixed text entries) so that
# values such as '43.40' are not silently normalised to the number 43.4.
$apos = [string][char]39

$ws.Range('D2').Value = '28.452.76'
$ws.Range('E2').Value = '  +5.36%  '
$ws.Range('D3').Value = '1.819.72'
$ws.Range('E3').Value = '  +5.35%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').Value = $apos + '318.34'
$ws.Range('E5').Value = '  +2.51%  '
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').Value = $apos + '0.5739'
$ws.Range('E7').Value = '  +18.73%  '
$ws.Range('D8').Value = $apos + '0.3857'
$ws.Range('E8').Value = '  +11.10%  '
$ws.Range('D9').Value = $apos + '43.40'
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('D10').Value = $apos + '0.07647'
$ws.Range('E10').Value = '  +5.78%  '
$ws.Range('E11').Value = '  +8.57%  '
$ws.Range('D12').Value = $apos + '21.41'
$ws.Range('E12').Value = '  +7.46%  '
$ws.Range('D13').Value = $apos + '0.9976'
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('E14').Value = '  +6.69%  '
$ws.Range('D15').Value = '1.813.84'
$ws.Range('E15').Value = '  +3.95%  '
$ws.Range('D16').Value = $apos + '7.307'
$ws.Range('E16').Value = '  +7.26%  '
$ws.Range('D17').Value = $apos + '92.42'
$ws.Range('E17').Value = '  +6.24%  '
$ws.Range('D18').Value = $apos + '0.00001084'
$ws.Range('E18').Value = '  +4.99%  '
$ws.Range('D19').Value = $apos + '0.06517'
$ws.Range('D20').Value = $apos + '0.9978'
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('E21').Value = '  +4.78%  '
$ws.Range('D22').Value = $apos + '6.013'
$ws.Range('E22').Value = '  +5.31%  '
$ws.Range('D23').Value = '28.465.64'
$ws.Range('E23').Value = '  +5.15%  '
$ws.Range('D24').Value = $apos + '11.38'
$ws.Range('E24').Value = '  +3.92%  '
$ws.Range('D25').Value = $apos + '2.093'
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('D26').Value = $apos + '21.01'
$ws.Range('E26').Value = '  +5.53%  '
$ws.Range('D27').Value = $apos + '157.64'
$ws.Range('E27').Value = '  +2.33%  '
$ws.Range('D28').Value = $apos + '2.408'
$ws.Range('E28').Value = '  +16.61%  '
$ws.Range('D29').Value = '2.022.44'
$ws.Range('E29').Value = '  +5.21%  '
$ws.Range('D30').Value = $apos + '123.85'
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('D31').Value = $apos + '1.168'
$ws.Range('E31').Value = '  +12.72%  '
$ws.Range('D32').Value = $apos + '0.1057'
$ws.Range('E32').Value = '  +13.56%  '
$ws.Range('E33').Value = '  +7.91%  '
$ws.Range('D34').Value = $apos + '3.634'
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('D35').Value = $apos + '0.02324'
$ws.Range('E35').Value = '  +6.78%  '
$ws.Range('D36').Value = $apos + '0.2170'
$ws.Range('E36').Value = '  +8.68%  '
$ws.Range('D37').Value = $apos + '8.806'
$ws.Range('E37').Value = '  +17.07%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = $apos + '11.76'
$ws.Range('E38').Value = '  +7.45%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = $apos + '0.6489'
$ws.Range('E39').Value = '  +8.65%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = $apos + '0.06105'
$ws.Range('E40').Value = '  +3.43%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').Value = $apos + '5.067'
$ws.Range('E41').Value = '  +6.89%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = $apos + '1.159'
$ws.Range('E42').Value = '  +3.43%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').Value = $apos + '0.9976'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').Value = $apos + '1.378'
$ws.Range('E44').Value = '  -3.65%  '
$ws.Range('D45').Value = $apos + '13.54'
$ws.Range('E45').Value = '  +5.94%  '
$ws.Range('D46').Value = $apos + '0.6040'
$ws.Range('E46').Value = '  +7.67%  '
$ws.Range('D47').Value = $apos + '3.723'
$ws.Range('E47').Value = '  +3.99%  '
$ws.Range('D48').Value = $apos + '122.51'
$ws.Range('E48').Value = '  +2.75%  '
$ws.Range('D49').Value = $apos + '1.955'
$ws.Range('E49').Value = '  +5.89%  '
$ws.Range('E50').Value = '  +4.49%  '
$ws.Range('D51').Value = $apos + '0.06858'
$ws.Range('E51').Value = '  +3.09%  '
